$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect to make the edits, then re-protect.
$ws.Unprotect()

# Update the confidential notice date from 2021-07-08 to 2021-07-09
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."

# Update the Weight/Percent Change figures for EFA, EEM, Total rows
$ws.Range("D2").Value = 0.8478623709366262
$ws.Range("E2").Value = 0.01700115045379014

$ws.Range("D3").Value = 0.1521376290633738
$ws.Range("E3").Value = 0.01728723404255317

$ws.Range("E4").Value = 0.01704467453269864

$ws.Protect()
